# Add a new "stock_data" worksheet with food-stock numbers, mirroring the
# class example workbook update (dataframe merging / concatenating demo).

$wb = $excel.ActiveWorkbook

# --- create the new sheet, placed after the existing "food_data" sheet ---
$wsFood = $wb.Worksheets.Item("food_data")
$wsStock = $wb.Worksheets.Add($null, $wsFood)
$wsStock.Name = "stock_data"

# --- fill in the stock data ---
$wsStock.Range("A1").Value = "Food"
$wsStock.Range("B1").Value = "Amount Available (lb)"

$wsStock.Range("A2").Value = "banana"
$wsStock.Range("B2").Value = 30

$wsStock.Range("A3").Value = "apple"
$wsStock.Range("B3").Value = 80

$wsStock.Range("A4").Value = "orange"
$wsStock.Range("B4").Value = 60

$wsStock.Range("A5").Value = "bell pepper"
$wsStock.Range("B5").Value = 20

$wsStock.Range("A6").Value = "kobe beef"
$wsStock.Range("B6").Value = 4

# --- selections / active sheet, matching the saved view state ---
$wsFood.Range("A1:A6").Select() | Out-Null
$wsStock.Activate() | Out-Null
$wsStock.Range("B6").Select() | Out-Null
